# Apply the "fixed workflow" re-run: the sensitivity sweep now starts at
# cutoff step 5 instead of step 1, so the first four rows of results
# (old cutoff steps 1-4) are dropped from each sheet and everything
# shifts up; the Cutoff index column (A) is renumbered sequentially.
$wb = $excel.ActiveWorkbook

ForEach ($ws in $wb.Worksheets) {
    # Drop the first four data rows (old rows 2:5) - shifts remaining
    # data rows up so row 2 now holds what used to be row 6, etc.
    $ws.Rows("2:5").Delete()

    # Renumber the Cutoff column (A) back to a clean 0..N-1 sequence
    # for the remaining 15 data rows (rows 2-16).
    For ($i = 0; $i -lt 15; $i++) {
        $ws.Cells.Item($i + 2, 1).Value = $i
    }
}
